$d = $word.ActiveDocument

function Set-ParagraphRuns {
    param($ParaIndex, $StyleId, $Tokens)

    $runsXml = ""
    foreach ($tok in $Tokens) {
        $escaped = $tok.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="' + $StyleId + '"/></w:pPr>' + $runsXml + '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $range = $d.Paragraphs($ParaIndex).Range
    $range.InsertXML($xml)
}

# Title: "Answers: Laws of indices" -> split into word/space runs
Set-ParagraphRuns 1 "Title" @(
    "Answers:", " ", "Laws", " ", "of", " ", "indices"
)

# Author: "Isabella Lewis, Akshat Srivastava" -> split into word/space runs
Set-ParagraphRuns 2 "Author" @(
    "Isabella", " ", "Lewis,", " ", "Akshat", " ", "Srivastava"
)

# Abstract: "Answers to questions relating to using laws of indices." -> split into word/space runs
Set-ParagraphRuns 4 "Abstract" @(
    "Answers", " ", "to", " ", "questions", " ", "relating", " ", "to", " ",
    "using", " ", "laws", " ", "of", " ", "indices."
)
